$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...around Personal Data Ecosystems [" becomes
#           "...around Personal Data Ecosystems and MyData [" with "MyData"
#           italicised.
# ---------------------------------------------------------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Personal Data Ecosystems [",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Personal Data Ecosystems and MyData [", 2)
if (-not $found1) {
    throw "Could not find 'Personal Data Ecosystems [' to replace"
}

# Now italicise just the newly inserted word "MyData" within that same
# sentence (search again so we get a fresh, correctly bounded Range).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Personal Data Ecosystems and MyData [",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
if (-not $found2) {
    throw "Could not find the freshly inserted 'MyData' text"
}
# Narrow the range down to just the word "MyData" inside the match.
$matchText = $rng2.Text
$offset = $matchText.IndexOf("MyData")
$myDataStart = $rng2.Start + $offset
$myDataEnd = $myDataStart + 6
$myDataRange = $d.Range($myDataStart, $myDataEnd)
if ($myDataRange.Text -ne "MyData") {
    throw "Range offset computation for 'MyData' is wrong: got '$($myDataRange.Text)'"
}
$myDataRange.Italic = 1

# ---------------------------------------------------------------------------
# Change 2: ". Collectively, I now knew where to position my existing and
#            newly discovered understandings..." becomes
#            ". Collectively through these discoveries, I learned knew where
#            to position my existing and newly-discovered understandings..."
# ---------------------------------------------------------------------------

$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "Collectively, I now knew where to position my existing and newly discovered understandings into the existing research landscape.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Collectively through these discoveries, I learned knew where to position my existing and newly-discovered understandings into the existing research landscape.",
    2)
if (-not $found3) {
    throw "Could not find the 'Collectively, I now knew...' sentence to replace"
}

Write-Output "done"
